$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.271.02'
$ws.Range("E2").Value = '  -2.02%  '
$ws.Range("D3").Value = '3.381.39'
$ws.Range("E3").Value = '  -1.80%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'567.62"
$ws.Range("E5").Value = '  -2.01%  '
$ws.Range("D6").Value = "'140.24"
$ws.Range("E6").Value = '  -6.43%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.380.90'
$ws.Range("E8").Value = '  -1.84%  '
$ws.Range("E9").Value = '  -0.62%  '
$ws.Range("E10").Value = '  -4.77%  '
$ws.Range("E11").Value = '  -1.48%  '
$ws.Range("E12").Value = '  -1.17%  '
$ws.Range("D13").Value = '3.957.89'
$ws.Range("E13").Value = '  -1.82%  '
$ws.Range("D14").Value = "'0.124"
$ws.Range("E14").Value = '  +1.21%  '
$ws.Range("D15").Value = "'27.96"
$ws.Range("E15").Value = '  -0.15%  '
$ws.Range("D16").Value = '3.385.66'
$ws.Range("E16").Value = '  -1.72%  '
$ws.Range("E17").Value = '  -3.48%  '
$ws.Range("D18").Value = '60.392.28'
$ws.Range("E18").Value = '  -1.96%  '
$ws.Range("D19").Value = "'6.21"
$ws.Range("E19").Value = '  -1.72%  '
$ws.Range("D20").Value = "'13.96"
$ws.Range("E20").Value = '  -2.28%  '
$ws.Range("D21").Value = "'9.02"
$ws.Range("E21").Value = '  -4.98%  '
$ws.Range("D22").Value = "'386.49"
$ws.Range("E22").Value = '  -0.45%  '
$ws.Range("E23").Value = '  -2.16%  '
$ws.Range("D24").Value = "'73.25"
$ws.Range("E24").Value = '  +0.48%  '
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").Value = "'0.0000116"
$ws.Range("E26").Value = '  -6.16%  '
$ws.Range("D27").Value = '3.529.70'
$ws.Range("E27").Value = '  -1.60%  '
$ws.Range("E28").Value = '  -1.58%  '
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  -5.16%  '
$ws.Range("D31").Value = "'7.89"
$ws.Range("E31").Value = '  -4.46%  '
$ws.Range("E32").Value = '  -2.22%  '
$ws.Range("E33").Value = '  -7.78%  '
$ws.Range("D35").Value = "'23.59"
$ws.Range("E35").Value = '  -1.78%  '
$ws.Range("D36").Value = '3.411.99'
$ws.Range("E36").Value = '  -1.64%  '
$ws.Range("D37").Value = "'6.89"
$ws.Range("E37").Value = '  -2.58%  '
$ws.Range("D38").Value = "'167.63"
$ws.Range("E38").Value = '  +0.64%  '
$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").Value = "'4.91"
$ws.Range("E39").Value = '  -7.65%  '
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").Value = "'1.49"
$ws.Range("E40").Value = '  -5.15%  '
$ws.Range("D41").Value = "'0.0770"
$ws.Range("E41").Value = '  -2.81%  '
$ws.Range("D42").Value = "'26.98"
$ws.Range("E42").Value = '  +1.93%  '
$ws.Range("E43").Value = '  -1.56%  '
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("E45").Value = '  -1.80%  '
$ws.Range("E46").Value = '  -1.84%  '
$ws.Range("D47").Value = "'41.20"
$ws.Range("E47").Value = '  -2.60%  '
$ws.Range("D48").Value = '2.511.70'
$ws.Range("E48").Value = '  -3.82%  '
$ws.Range("D49").Value = "'1.12"
$ws.Range("E49").Value = '  -3.91%  '
$ws.Range("E50").Value = '  -3.95%  '
$ws.Range("D51").Value = "'22.92"
$ws.Range("E51").Value = '  -1.61%  '
